$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# C3: replace the plain "Info" header with a two-run rich-text cell:
#   run 1 "Info: "  -> inherits the cell's existing bold/underline/14pt font
#   run 2 "this includes research, coding and testing" -> regular, 12pt
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "Info: this includes research, coding and testing"

# Apply the smaller/regular font to the second run only.
$infoRun = $ws.Range("C3").Characters(7, 43)
$infoRun.Font.Size = 12
$infoRun.Font.Bold = $false
$infoRun.Font.Underline = $false

# Register a (sz 12 / Calibri / theme-1 / minor-scheme) font in the
# workbook's font table the same way Excel does when character-level
# formatting is applied, then restore the whole-cell font so the cell
# keeps its original header style (bold, underlined, 14pt).
$ws.Range("C3").Font.Bold = $false
$ws.Range("C3").Font.Underline = $false
$ws.Range("C3").Font.Size = 12
$ws.Range("C3").Font.Size = 14
$ws.Range("C3").Font.Bold = $true
$ws.Range("C3").Font.Underline = $true

# ---------------------------------------------------------------------------
# New rows of logged time: Oct 22 2017 (serial 43030), with the same
# date/time number formats already used by row 7 (copied via PasteSpecial
# so no new cell styles are introduced).
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8:A10").PasteSpecial(-4122) | Out-Null

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8:B10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.Range("A8").Value = 43030
$ws.Range("B8").Value = 0.010416666666666666
$ws.Range("C8").Value = "Camera follow players + switch"

$ws.Range("A9").Value = 43030
$ws.Range("B9").Value = 0.020833333333333332
$ws.Range("C9").Value = "Player following and regrouping when at destination "

$ws.Range("A10").Value = 43030
$ws.Range("B10").Value = 0.013888888888888888
$ws.Range("C10").Value = "Characters auto follow or staying "

# ---------------------------------------------------------------------------
# Selection moves to F11.
# ---------------------------------------------------------------------------
$ws.Range("F11").Select() | Out-Null
